$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A121").Value = "Al Siraat College, EPPINGCharles La Trobe, MACLEOD WESTEpping Secondary College, EPPINGEpping Views Primary School, EPPINGFitzroy Primary School, FITZROYGisborne Secondary College, GISBORNEGladstone Park Secondary College, GLADSTONE PARKGreenvalePrimary School, GREENVALELalor Secondary College, LALORMelbourne Girls College, RICHMONDNorthcote High School, NORTHCOTENorthern College of the Arts and Technology, PRESTONPascoe Vale Girls Secondary College, OAK PARKPenders Grove Primary School, THORNBURYPeter Lalor Secondary College, LALORPrinces Hill Secondary College, PRINCES HILLRoxburgh College, ROXBURGH PARKSacred Heart School, FITZROYSt Joseph's Primary School, MERNDAThornbury High School, THORNBURYTAFE"

$ws.Range("A174").Value = "Al Taqwa College, TRUGANINABelmont High School, BELMONTBraybrook College, BRAYBROOKBuckley Park College, ESSENDONCana Catholic Primary School, HILLSIDECaroline Chilsholm Catholic College, BRAYBROOKCatholic Regional College, SYDENHAMClonard Secondary College, GEELONG WESTColac Secondary College, COLACCopperfield College Delahey Campus, DELAHEYCopperfield College, Sydenham Campus, SYDENHAM"

$ws.Range("A177").Value = "Grovedale West Primary School, GROVEDALEHoppers Crossing Secondary College, HOPPERS CROSSINGLowther Hall Anglican Grammar School, ESSENDONKeilor Downs Secondary College, KEILOR DOWNSKeilor Primary School, KEILORMarian College, SUNSHINE WESTManor Lakes P-12 College, WYNDHAM VALEMacKillop College, WERRIBEEMelton Secondary College, MELTONMount St Joseph's Girls' College, ALTONAPoint Cook Senior Secondary College, POINT COOKRosamond Specialist School, BRAYBROOKSt Alban's Secondary College, ST ALBANSSt Joseph's Flexible Learning Centre Geelong Campus, GEELONGStaughton College, MELTON SOUTHSunshine Heights Primary School, SUNSHINESydenham Hillside Primary School, SYDENHAMTarneit Rise Primary School, TARNEITTarneit Senior College, TARNEITTaylors Lakes Secondary College, TAYLORS LAKESThomas Carr College, TARNEITWarracknabeal Secondary College, WARRACKNABEALTAFE"

$ws.Range("A187").Value = "li: Last Update: 30 July 2020"
